# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (interest counts) and sold-out status text
# across the "展览" and "全部类型" sheets (and matching rows in
# "本地生活"), matching the source data refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 1146
$ws.Range("F7").Value  = 805
$ws.Range("F8").Value  = 271
$ws.Range("F10").Value = 676
$ws.Range("F11").Value = 432
$ws.Range("F14").Value = 953
$ws.Range("F15").Value = 132
$ws.Range("F16").Value = 2070
$ws.Range("F17").Value = 547
$ws.Range("F18").Value = 9177
$ws.Range("F19").Value = 855
$ws.Range("F20").Value = 533
$ws.Range("F21").Value = 83
$ws.Range("F22").Value = 107

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5626
$ws.Range("G2").Value = "已售罄"
$ws.Range("F3").Value = 436
$ws.Range("F4").Value = 414

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 5626
$ws.Range("G3").Value  = "已售罄"
$ws.Range("F4").Value  = 436
$ws.Range("F5").Value  = 414
$ws.Range("F7").Value  = 1146
$ws.Range("F12").Value = 805
$ws.Range("F14").Value = 271
$ws.Range("F16").Value = 682
$ws.Range("F17").Value = 432
$ws.Range("F22").Value = 953
$ws.Range("F24").Value = 132
$ws.Range("F27").Value = 2070
$ws.Range("F28").Value = 547
$ws.Range("F29").Value = 9177
$ws.Range("F32").Value = 855
$ws.Range("F33").Value = 533
$ws.Range("F34").Value = 83
$ws.Range("F35").Value = 107

$wb.Save()
